# feat: Added ship api call
#
# Updates the "ships" sheet: currentFuel (D), maxFuel (E) and
# fuelRechargeRate (F) go from 5 -> 25 for every ship row (rows 2-31),
# then leaves the "ships" sheet active/selected (it becomes the 4th tab,
# selection sits on F32).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ships")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("D$r").Value = 25
    $ws.Range("E$r").Value = 25
    $ws.Range("F$r").Value = 25
}

$ws.Activate()
$ws.Range("F32").Select()
